$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

function Split-NewRun($range) {
    # Forces a run boundary at the edges of $range: toggling Bold on then off
    # is a formatting no-op but makes the engine split the containing run(s)
    # so that $range becomes (or stays) its own <w:r>.
    $range.Font.Bold = 1
    $range.Font.Bold = 0
}

# Rebuilds the run layout of a paragraph so it contains exactly the given
# ordered list of text segments (each a plain string; "`v" i.e. chr(11)
# inside a segment is left alone -- there are none here) as separate runs,
# without disturbing any <w:br/> runs that live between them.
#
# $paraIndex   : 1-based Paragraphs index of the target paragraph
# $segments    : ordered array of the exact text runs the paragraph's
#                visible text should be split into once <w:br/> chars are
#                excluded (concatenating them with the existing break
#                characters reproduces Paragraph.Range.Text exactly)
function Rebuild-ParagraphRuns($paraIndex, [string[]]$segments) {
    $para = $d.Paragraphs.Item($paraIndex)
    $pStart = $para.Range.Start
    $fullText = $para.Range.Text

    $offsets = @()
    $searchFrom = 0
    foreach ($seg in $segments) {
        $idx = $fullText.IndexOf($seg, $searchFrom)
        if ($idx -lt 0) {
            throw "segment not found: $seg"
        }
        $segStart = $idx
        $segEnd = $idx + $seg.Length
        $offsets += ,@($segStart, $segEnd)
        $searchFrom = $segEnd
    }

    for ($i = $offsets.Length - 1; $i -ge 0; $i--) {
        $pair = $offsets[$i]
        if ($pair[1] -gt $pair[0]) {
            $r = $d.Range($pStart + $pair[0], $pStart + $pair[1])
            Split-NewRun $r
        }
    }
}

# ============================================================
# Title / byline / contact info
# ============================================================
Replace-Text "Celestial Symphony: The Rhythms of the Universe" "The Allure of History: A Journey Through Time"
Replace-Text "Isabella Maxwell" "Alex Watson"
Replace-Text "isabella" "alex"
Replace-Text "maxwell@astronomycenter" "watson@eduworld"

# ============================================================
# Body paragraph -- simple 1-for-1 sentence replacements
# ============================================================
Replace-Text "Within the vast canvas of the cosmic tapestry, celestial bodies engage in an intricate dance, governed by the harmonious laws of gravitation" "History beckons us, like an alluring whisper from times gone by"
Replace-Text " From the gentle waltz of our solar system's planets to the whirling dervishes of distant galaxies, the cosmos pulsates with an unseen symphony" " It is a kaleidoscope of human experiences, triumphs and follies, wisdom and folly, painted on the canvas of centuries"
Replace-Text " This symphony is a symphony of motion, dictated by the gravitational forces that bind celestial bodies together" " With each chapter, it holds a mirror to our present, shedding light on our origins, shaping our identities, and guiding us into the future"
Replace-Text "It is a symphony of time, measured by the ebb and flow of stellar ages, the rise and fall of civilizations, and the birth and death of stars" "In its vast expanse, history unveils a tapestry woven with countless threads, each representing the lives of individuals who have shaped our world"
Replace-Text " Every celestial object, from the smallest comet to the grandest supermassive black hole, plays a role in this cosmic orchestra, contributing to the intricate composition that orchestrates the universe" " Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures"
Replace-Text "The study of these celestial motions, known as celestial mechanics, delves into the underlying principles that govern the dynamics of the universe" "Furthermore, history teaches us the art of empathy and perspective"
Replace-Text " Scientists, like maestros of the universe, analyze the ballet of planets, the pirouette of stars, and the majestic procession of galaxies" " As we journey through the annals of time, we encounter diverse cultures, beliefs, and ways of life"
Replace-Text " Through this meticulous examination, they unravel the mysteries of the cosmos, revealing its hidden harmonies and unlocking its secrets" " We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do"

# ============================================================
# Summary paragraph -- simple 1-for-1 sentence replacements
# ============================================================
Replace-Text "The symphony of the universe unfolds through the graceful dance of celestial bodies, governed by the gravitational forces that orchestrate the cosmos" "History is a tapestry of human experiences, unveiling the richness and complexity of our shared past"
Replace-Text " Celestial mechanics, like a conductor of the universal orchestra, analyzes this dance, deciphering the principles that govern its rhythm" " It illuminates our present, shaping our identities, and guiding us into the future"
Replace-Text " The study of this cosmic choreography grants us insights into the workings of the universe, allowing us to appreciate the profound beauty and intricate interconnectedness of all things celestial" " Through its stories of triumphs and follies, wisdom and folly, history teaches us empathy, perspective, and the profound interconnectedness of humanity"

# ============================================================
# New sentences inserted into the body paragraph
# ============================================================
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures") | Out-Null
$anchor.Collapse(0)
$anchor.Text = "."
$anchor.Collapse(0)
$anchor.Text = " Amidst the ebb and flow of civilizations and empires, we discover the timeless struggles of humanity - the quest for power, the pursuit of justice, the yearning for freedom"

$anchor = $d.Content.Duplicate
$anchor.Find.Execute("We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do") | Out-Null
$anchor.Collapse(0)
$anchor.Text = "."
$anchor.Collapse(0)
$anchor.Text = " This understanding fosters tolerance, compassion, and the realization that we are all part of a shared human story"

# ============================================================
# New sentences inserted into the summary paragraph
# ============================================================
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Through its stories of triumphs and follies, wisdom and folly, history teaches us empathy, perspective, and the profound interconnectedness of humanity") | Out-Null
$anchor.Collapse(0)
$anchor.Text = "."
$anchor.Collapse(0)
$anchor.Text = " It is a beacon that sheds light on our origins, offering invaluable insights into who we are and how we can navigate the challenges of an ever-changing world"

# ============================================================
# Rebuild run boundaries precisely for both edited paragraphs
# (text content is already final at this point; this just makes
#  sure every sentence/period lives in its own <w:r>).
# ============================================================
$bodyParaIndex = 5
Rebuild-ParagraphRuns $bodyParaIndex @(
    "History beckons us, like an alluring whisper from times gone by",
    ".",
    " It is a kaleidoscope of human experiences, triumphs and follies, wisdom and folly, painted on the canvas of centuries",
    ".",
    " With each chapter, it holds a mirror to our present, shedding light on our origins, shaping our identities, and guiding us into the future",
    ".",
    "In its vast expanse, history unveils a tapestry woven with countless threads, each representing the lives of individuals who have shaped our world",
    ".",
    " Like intricate brushstrokes, their actions, decisions, and sacrifices add color and texture, creating a vibrant masterpiece that reveals the human capacity for both great achievements and devastating failures",
    ".",
    " Amidst the ebb and flow of civilizations and empires, we discover the timeless struggles of humanity - the quest for power, the pursuit of justice, the yearning for freedom",
    ".",
    "Furthermore, history teaches us the art of empathy and perspective",
    ".",
    " As we journey through the annals of time, we encounter diverse cultures, beliefs, and ways of life",
    ".",
    " We learn to appreciate the richness of human existence, and we begin to understand why people think, feel, and act as they do",
    ".",
    " This understanding fosters tolerance, compassion, and the realization that we are all part of a shared human story",
    "."
)

$summaryParaIndex = 7
Rebuild-ParagraphRuns $summaryParaIndex @(
    "History is a tapestry of human experiences, unveiling the richness and complexity of our shared past",
    ".",
    " It illuminates our present, shaping our identities, and guiding us into the future",
    ".",
    " Through its stories of triumphs and follies, wisdom and folly, history teaches us empathy, perspective, and the profound interconnectedness of humanity",
    ".",
    " It is a beacon that sheds light on our origins, offering invaluable insights into who we are and how we can navigate the challenges of an ever-changing world",
    "."
)

# ============================================================
# New empty paragraph at the end of the document body
# ============================================================
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

Write-Output "edit complete"
